$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 19, shifting existing rows 19:78 down to 20:79
$ws.Rows.Item(19).Insert()

# Populate the newly inserted row 19 with the new weekly record
$ws.Range("A19").Value = 3
$ws.Range("B19").Value = "Femacal de La Calera"
$ws.Range("C19").Value = "Coquimbo"
$ws.Range("D19").Value = (Get-Date -Year 2022 -Month 12 -Day 15 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$ws.Range("E19").Value = 5
$ws.Range("F19").Value = 100112022
$ws.Range("G19").Value = "Arveja Verde"
$ws.Range("H19").Value = "Perfection"
$ws.Range("I19").Value = "Primera"
$ws.Range("J19").Value = 35
$ws.Range("K19").Value = 19000
$ws.Range("L19").Value = 19000
$ws.Range("M19").Value = 19000
$ws.Range("N19").Value = "`$/malla 25 kilos"
$ws.Range("O19").Value = "Provincia de Limarí"
$ws.Range("P19").Value = 760
$ws.Range("Q19").Value = 25
$ws.Range("R19").Value = "Hortaliza"
